# Update results for row 380 kV case (Case_5_232, pl_mw sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new values for columns B..O (columns 2..15)
$newValues = @{
    2 = @(0.4002345243975469, 0.0292379035167869, 0.1684952640472517, 0.1558357427434238, 1.589051218795753, 0, 0.07973214163530429, 0, 0.1840581092676743, 0.3524572108632924, 0, 0.2310193804041845, 1.941613291519175, 3.949852420747931)
    3 = @(0.3678260096117185, 0.02554740063331451, 0.1653073817920472, 0.1546393167853886, 1.592261207634152, 0, 0.07973214163530429, 0, 0.1838083673128992, 0.3189940532914477, 0, 0.2208632601734948, 1.960455524992266, 3.966629397327495)
    4 = @(0.348042742285628, 0.02327198865151558, 0.1634166907763728, 0.1539715075597314, 1.594967748583379, 0, 0.07973214163530429, 0, 0.1837366976301169, 0.2985130855275884, 0, 0.2147240888939876, 1.972615078615114, 3.978778437857443)
    5 = @(0.3400104637549646, 0.02234241728579889, 0.1626630749418041, 0.1537162120047846, 1.596255840223371, 0, 0.07973214163530429, 0, 0.1837280565285297, 0.290183874743775, 0, 0.2122468196736591, 1.977718665749172, 3.98419421570226)
    6 = @(0.3386785094795357, 0.02218792381748358, 0.162538958276059, 0.1536748389022264, 1.596480915109105, 0, 0.07973214163530429, 0, 0.1837278645598772, 0.2888018531401286, 0, 0.2118369553211288, 1.978575080459729, 3.985121591112957)
    7 = @(0.3479342958412701, 0.02325946146702051, 0.1634064588819513, 0.1539679963022955, 1.594984370313135, 0, 0.07973214163530429, 0, 0.183736497790008, 0.2984006856286499, 0, 0.2146905801785834, 1.972683306195194, 3.978849594025561)
    8 = @(0.3890363294037229, 0.02796740286039778, 0.1673822942559724, 0.1554093764511997, 1.590005457046189, 0, 0.07973214163530429, 0, 0.1839550630701439, 0.3409057801354152, 0, 0.2274975727689963, 1.947987582870733, 3.955253747856034)
    9 = @(0.4705386873024793, 0.0371231546221793, 0.1757046279797692, 0.1587643007657995, 1.5860717657782, 0, 0.07973214163530429, 0, 0.1850307187896547, 0.4247618543171257, 0, 0.2533735112605555, 1.904241703089423, 3.923635116581323)
    10 = @(0.5309518586561524, 0.04380158918949917, 0.1821358201218288, 0.1615494938843902, 1.586728821723767, 0, 0.07973214163530429, 0, 0.1862144641844665, 0.486662157678353, 0, 0.2728428026374203, 1.874952586932343, 3.909328351868169)
    11 = @(0.5585481322935379, 0.04682898100824673, 0.185129576148455, 0.1628857550421543, 1.587796632561094, 0, 0.07973214163530429, 0, 0.1868382272286908, 0.5148824520858852, 0, 0.2817982479464689, 1.862246525303699, 3.904755804067605)
    12 = @(0.5690141187568543, 0.04797380137172524, 0.186272957466258, 0.1634016835525642, 1.588311408281868, 0, 0.07973214163530429, 0, 0.1870866697284299, 0.5255771964203859, 0, 0.2852035000611366, 1.857523878761645, 3.903302457282138)
    13 = @(0.5667593834725437, 0.04772731527637575, 0.1860262794221796, 0.1632901288357935, 1.588195633595504, 0, 0.07973214163530429, 0, 0.1870326193459704, 0.5232735295836335, 0, 0.2844694974983142, 1.858537032341631, 3.903603092176439)
    14 = @(0.5594088600857106, 0.04692319818893509, 0.1852234487276405, 0.1629280023268187, 1.587836771368302, 0, 0.07973214163530429, 0, 0.1868584216385543, 0.5157621511338562, 0, 0.2820781201698424, 1.861856209151061, 3.904630662224463)
    15 = @(0.5549085044420679, 0.04643044494777371, 0.1847329536412161, 0.1627074793407814, 1.587631333247501, 0, 0.07973214163530429, 0, 0.186753313386987, 0.5111622845591342, 0, 0.2806151527934233, 1.863900874579154, 3.905296300491386)
    16 = @(0.5291506272585593, 0.04360352285107183, 0.1819415350629328, 0.1614635560184468, 1.586674497465054, 0, 0.07973214163530429, 0, 0.1861754135485754, 0.4848190891079014, 0, 0.2722595150138858, 1.87579539092807, 3.909666114794078)
    17 = @(0.5133778376961118, 0.04186652852803263, 0.1802464854936829, 0.1607181575898906, 1.586284340195746, 0, 0.07973214163530429, 0, 0.1858427130199232, 0.4686737847660538, 0, 0.267158768205924, 1.883250526598499, 3.912842505613838)
    18 = @(0.5043165094056405, 0.04086645505074671, 0.1792779602543959, 0.1602959461419147, 1.586132314727152, 0, 0.07973214163530429, 0, 0.1856593801112041, 0.4593932557598066, 0, 0.2642342609464663, 1.887596662243897, 3.914851686231628)
    19 = @(0.5012503637226473, 0.04052767720452266, 0.1789511401283903, 0.1601541141164979, 1.586093277872948, 0, 0.07973214163530429, 0, 0.1855986862519146, 0.4562520439689308, 0, 0.2632456771700546, 1.889078173374548, 3.915563259258107)
    20 = @(0.5150557678902032, 0.04205153845622078, 0.18042626238541, 0.1607968318436619, 1.586318382851516, 0, 0.07973214163530429, 0, 0.1858772989096451, 0.4703918810612606, 0, 0.2677007894574146, 1.882450897262995, 3.912485517116352)
    21 = @(0.561567460241605, 0.04715943037780335, 0.1854589967511515, 0.1630340989143164, 1.587939182360145, 0, 0.07973214163530429, 0, 0.1869092558191099, 0.5179682030808976, 0, 0.2827801466404551, 1.860878874424416, 3.904321291675046)
    22 = @(0.5920578391312574, 0.0504884564207515, 0.1888047390514629, 0.1645540561624514, 1.589642004917593, 0, 0.07973214163530429, 0, 0.1876550141399207, 0.549110400445727, 0, 0.2927170083948312, 1.847298292990333, 3.900606853225526)
    23 = @(0.5757762764701511, 0.04871255924722107, 0.1870139088181588, 0.1637375547857687, 1.58867434239302, 0, 0.07973214163530429, 0, 0.1872504718321579, 0.5324849734143697, 0, 0.2874061097376028, 1.85449909617367, 3.902441020758602)
    24 = @(0.5142971548705759, 0.0419678999893307, 0.1803449666153369, 0.1607612434994223, 1.586302767004256, 0, 0.07973214163530429, 0, 0.1858616378916977, 0.4696151246366753, 0, 0.2674557167794163, 1.88281222254334, 3.912646341469269)
    25 = @(0.448395304279984, 0.03465464304244392, 0.1733973179343877, 0.1578003222436593, 1.58651251680034, 0, 0.07973214163530429, 0, 0.1846705573212759, 0.402024273450337, 0, 0.2462924803521744, 1.915575217445703, 3.930621077509926)
}

foreach ($rowNum in $newValues.Keys) {
    $values = $newValues[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # Column B is index 2
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}
